$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Refresh the recalculated capital-cost figures (rows 3-17, column C)
# ---------------------------------------------------------------------
$ws.Range("C3").Value  = 44.371165973952337
$ws.Range("C4").Value  = 26.53598744217307
$ws.Range("C5").Value  = 1.774846638958093
$ws.Range("C6").Value  = 3.9934049376557099
$ws.Range("C7").Value  = 1.996702468827855
$ws.Range("C8").Value  = 78.672107461567066
$ws.Range("C9").Value  = 7.8672107461567071
$ws.Range("C10").Value = 7.8672107461567071
$ws.Range("C11").Value = 15.734421492313411
$ws.Range("C12").Value = 23.60163223847012
$ws.Range("C13").Value = 7.8672107461567071
$ws.Range("C14").Value = 62.937685969253643
$ws.Range("C15").Value = 141.6097934308207
$ws.Range("C16").Value = 7.0804896715410353
$ws.Range("C17").Value = 148.69028310236169

# ---------------------------------------------------------------------
# 2. Rebuild the "Raw materials / By-products and credits" block
#    (rows 21-29) to add the new fatty-alcohols / Tridecane line and
#    re-order the existing raw materials.
# ---------------------------------------------------------------------

# The old layout has a single merged range A21:A26 - unmerge it first so
# every row in the block becomes independently addressable.
$ws.Range("A21:A26").UnMerge()

# --- Row 21: Raw materials / Tridecane -------------------------------
$ws.Range("A21").Value = "Raw materials"
$ws.Range("B21").Value = "Tridecane"
$ws.Range("C21").Value = 775.64317499999993
$ws.Range("D21").Value = -0.037139309833833033

# --- Row 22: (merged with 21) / Process water -------------------------
$ws.Range("A22").Value = ""
$ws.Range("B22").Value = "Process water"
$ws.Range("C22").Value = 0.320236305
$ws.Range("D22").Value = 0.3544937424716692

# --- Row 23: By-products and credits / Wastewater ----------------------
$ws.Range("A23").Value = "By-products and credits"
$ws.Range("B23").Value = "Wastewater"
$ws.Range("C23").Value = -3.4222048137973369
$ws.Range("D23").Value = -4.1054461435065672

# --- Row 24: Raw materials / Glucose -----------------------------------
$ws.Range("A24").Value = "Raw materials"
$ws.Range("B24").Value = "Glucose"
$ws.Range("C24").Value = 299.99700764999989
$ws.Range("D24").Value = 160.5470037981996

# --- Row 25: (merged with 24) / CSL -------------------------------------
$ws.Range("A25").Value = ""
$ws.Range("B25").Value = "CSL"
$ws.Range("C25").Value = 51.528108000000003
$ws.Range("D25").Value = 0.2436875313576996

# --- Row 26: (merged with 24) / DAP ---------------------------------------
$ws.Range("A26").Value = ""
$ws.Range("B26").Value = "DAP"
$ws.Range("C26").Value = 895.39159499999994
$ws.Range("D26").Value = 0.42241247900968393

# --- Row 27: (merged with 24) / Salt ---------------------------------------
$ws.Range("A27").Value = ""
$ws.Range("B27").Value = "Salt"
$ws.Range("C27").Value = 136.07775000000001
$ws.Range("D27").Value = 1.5067104172463091

# --- Row 28: (merged with 24) / Natural gas ---------------------------------
$ws.Range("A28").Value = ""
$ws.Range("B28").Value = "Natural gas"
$ws.Range("C28").Value = 197.76633000000001
$ws.Range("D28").Value = 3.0227578626756069

# --- Row 29: Total variable operating cost ------------------------------
$ws.Range("A29").Value = "Total variable operating cost"
$ws.Range("D29").Value = 170.16537266463331

# Re-merge the label column for each raw-materials group.
$ws.Range("A21:A22").Merge()
$ws.Range("A24:A28").Merge()

# Match the bold / bordered / centred-top look used by the rest of the
# table's label column for the two (now split) merged header cells.
$ws.Range("A21:A22,A24:A28").Font.Bold = $true
$ws.Range("A21:A22,A24:A28").HorizontalAlignment = -4108
$ws.Range("A21:A22,A24:A28").VerticalAlignment = -4160
$ws.Range("A21:A22,A24:A28").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3. Refresh maintenance / property-insurance figures (rows 35-36)
# ---------------------------------------------------------------------
$ws.Range("C35").Value = 1.33113497921857
$ws.Range("D35").Value = 1.2778895800498269
$ws.Range("C36").Value = 0.31059816181766642
$ws.Range("D36").Value = 0.29817423534495968
